# Update TPM-derived metrics for the Nampt-Insr LR-pair sheet.
# Columns G,H,I,J,M,N,O,P,Q,R,S,T for rows 2-26 are refreshed with
# values recomputed from the new TPM table (per commit "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "G2" = 9.508709
    "H2" = 28.526127
    "I2" = 0.05211107427675058
    "J2" = 0.05288783704580559
    "M2" = 15.70818033333333
    "N2" = 47.12454099999999
    "O2" = 0.3220467100482788
    "P2" = 0.334408980496766
    "Q2" = 149.3645157091896
    "R2" = 1344.280641382707
    "S2" = 0.01678220002790901
    "T2" = 0.01768616766716694
    "G3" = 9.508709
    "H3" = 28.526127
    "I3" = 0.05211107427675058
    "J3" = 0.05288783704580559
    "M3" = 5.85326
    "N3" = 17.55978
    "O3" = 0.1200026410479322
    "P3" = 0.1246091315254933
    "Q3" = 55.65694604133999
    "R3" = 500.91251437206
    "S3" = 0.006253466541055035
    "T3" = 0.006590307442539647
    "G4" = 9.508709
    "H4" = 28.526127
    "I4" = 0.05211107427675058
    "J4" = 0.05288783704580559
    "M4" = 10.959131
    "N4" = 32.877393
    "O4" = 0.2246824271585863
    "P4" = 0.2333072161810874
    "Q4" = 104.207187571879
    "R4" = 937.8646881469111
    "S4" = 0.01170844265034169
    "T4" = 0.01233911403099589
    "G5" = 9.508709
    "H5" = 28.526127
    "I5" = 0.05211107427675058
    "J5" = 0.05288783704580559
    "M5" = 5.4093935
    "N5" = 10.818787
    "O5" = 0.1109025579706895
    "P5" = 0.07677315161290731
    "Q5" = 51.4363486579915
    "R5" = 308.618091947949
    "S5" = 0.005779251435892236
    "T5" = 0.004060365931996368
    "G6" = 9.508709
    "H6" = 28.526127
    "I6" = 0.05211107427675058
    "J6" = 0.05288783704580559
    "M6" = 10.84612833333333
    "N6" = 32.538385
    "O6" = 0.2223656637745133
    "P6" = 0.230901520183746
    "Q6" = 103.1326780983217
    "R6" = 928.194102884895
    "S6" = 0.01158771362155261
    "T6" = 0.01221188197310675
    "G7" = 20.99754233333333
    "H7" = 62.992627
    "I7" = 0.1150739272977591
    "J7" = 0.1167892084285825
    "M7" = 15.70818033333333
    "N7" = 47.12454099999999
    "O7" = 0.3220467100482788
    "P7" = 0.334408980496766
    "Q7" = 329.8331815288007
    "R7" = 2968.498633759206
    "S7" = 0.03705917969857815
    "T7" = 0.03905536012362657
    "G8" = 20.99754233333333
    "H8" = 62.992627
    "I8" = 0.1150739272977591
    "J8" = 0.1167892084285825
    "M8" = 5.85326
    "N8" = 17.55978
    "O8" = 0.1200026410479322
    "P8" = 0.1246091315254933
    "Q8" = 122.9040746380066
    "R8" = 1106.13667174206
    "S8" = 0.01380917519148884
    "T8" = 0.01455300183383549
    "G9" = 20.99754233333333
    "H9" = 62.992627
    "I9" = 0.1150739272977591
    "J9" = 0.1167892084285825
    "M9" = 10.959131
    "N9" = 32.877393
    "O9" = 0.2246824271585863
    "P9" = 0.2333072161810874
    "Q9" = 230.1148171090457
    "R9" = 2071.033353981411
    "S9" = 0.02585508928793123
    "T9" = 0.02724776509846536
    "G10" = 20.99754233333333
    "H10" = 62.992627
    "I10" = 0.1150739272977591
    "J10" = 0.1167892084285825
    "M10" = 5.4093935
    "N10" = 10.818787
    "O10" = 0.1109025579706895
    "P10" = 0.07677315161290731
    "Q10" = 113.5839690139082
    "R10" = 681.503814083449
    "S10" = 0.01276199289305464
    "T10" = 0.008966275605438993
    "G11" = 20.99754233333333
    "H11" = 62.992627
    "I11" = 0.1150739272977591
    "J11" = 0.1167892084285825
    "M11" = 10.84612833333333
    "N11" = 32.538385
    "O11" = 0.2223656637745133
    "P11" = 0.230901520183746
    "Q11" = 227.7420388319327
    "R11" = 2049.678349487395
    "S11" = 0.02558849022670629
    "T11" = 0.02696680576721605
    "G12" = 76.01018666666667
    "H12" = 228.03056
    "I12" = 0.4165625936366697
    "J12" = 0.4227718364551835
    "M12" = 15.70818033333333
    "N12" = 47.12454099999999
    "O12" = 0.3220467100482788
    "P12" = 0.334408980496766
    "Q12" = 1193.981719330329
    "R12" = 10745.83547397296
    "S12" = 0.1341526128098675
    "T12" = 0.1413786988117234
    "G13" = 76.01018666666667
    "H13" = 228.03056
    "I13" = 0.4165625936366697
    "J13" = 0.4227718364551835
    "M13" = 5.85326
    "N13" = 17.55978
    "O13" = 0.1200026410479322
    "P13" = 0.1246091315254933
    "Q13" = 444.9073852085333
    "R13" = 4004.1664668768
    "S13" = 0.04998861139817693
    "T13" = 0.05268123137411832
    "G14" = 76.01018666666667
    "H14" = 228.03056
    "I14" = 0.4165625936366697
    "J14" = 0.4227718364551835
    "M14" = 10.959131
    "N14" = 32.877393
    "O14" = 0.2246824271585863
    "P14" = 0.2333072161810874
    "Q14" = 833.0055930144534
    "R14" = 7497.050337130081
    "S14" = 0.09359429460176283
    "T14" = 0.09863572024312482
    "G15" = 76.01018666666667
    "H15" = 228.03056
    "I15" = 0.4165625936366697
    "J15" = 0.4227718364551835
    "M15" = 5.4093935
    "N15" = 10.818787
    "O15" = 0.1109025579706895
    "P15" = 0.07677315161290731
    "Q15" = 411.1690096884533
    "R15" = 2467.01405813072
    "S15" = 0.04619785718921152
    "T15" = 0.03245752629784106
    "G16" = 76.01018666666667
    "H16" = 228.03056
    "I16" = 0.4165625936366697
    "J16" = 0.4227718364551835
    "M16" = 10.84612833333333
    "N16" = 32.538385
    "O16" = 0.2223656637745133
    "P16" = 0.230901520183746
    "Q16" = 824.4162392272889
    "R16" = 7419.7461530456
    "S16" = 0.09262921763765089
    "T16" = 0.09761865972837591
    "G17" = 8.0398025
    "H17" = 16.079605
    "I17" = 0.04406094930951247
    "J17" = 0.02981181178226265
    "M17" = 15.70818033333333
    "N17" = 47.12454099999999
    "O17" = 0.3220467100482788
    "P17" = 0.334408980496766
    "Q17" = 126.2906675143842
    "R17" = 757.7440050863049
    "S17" = 0.01418968376673247
    "T17" = 0.009969337584867932
    "G18" = 8.0398025
    "H18" = 16.079605
    "I18" = 0.04406094930951247
    "J18" = 0.02981181178226265
    "M18" = 5.85326
    "N18" = 17.55978
    "O18" = 0.1200026410479322
    "P18" = 0.1246091315254933
    "Q18" = 47.05905438115
    "R18" = 282.3543262869
    "S18" = 0.005287430284220563
    "T18" = 0.003714823975389219
    "G19" = 8.0398025
    "H19" = 16.079605
    "I19" = 0.04406094930951247
    "J19" = 0.02981181178226265
    "M19" = 10.959131
    "N19" = 32.877393
    "O19" = 0.2246824271585863
    "P19" = 0.2333072161810874
    "Q19" = 88.10924881162751
    "R19" = 528.6554928697651
    "S19" = 0.0098997210337727
    "T19" = 0.006955310816234241
    "G20" = 8.0398025
    "H20" = 16.079605
    "I20" = 0.04406094930951247
    "J20" = 0.02981181178226265
    "M20" = 5.4093935
    "N20" = 10.818787
    "O20" = 0.1109025579706895
    "P20" = 0.07677315161290731
    "Q20" = 43.49045538478375
    "R20" = 173.961821539135
    "S20" = 0.004886471985041817
    "T20" = 0.002288746745815107
    "G21" = 8.0398025
    "H21" = 16.079605
    "I21" = 0.04406094930951247
    "J21" = 0.02981181178226265
    "M21" = 10.84612833333333
    "N21" = 32.538385
    "O21" = 0.2223656637745133
    "P21" = 0.230901520183746
    "Q21" = 87.20072968965417
    "R21" = 523.204378137925
    "S21" = 0.009797642239744922
    "T21" = 0.006883592659956157
    "G22" = 67.91378400000001
    "H22" = 203.741352
    "I22" = 0.372191455479308
    "J22" = 0.3777393062881658
    "M22" = 15.70818033333333
    "N22" = 47.12454099999999
    "O22" = 0.3220467100482788
    "P22" = 0.334408980496766
    "Q22" = 1066.801966191048
    "R22" = 9601.217695719432
    "S22" = 0.1198630337451916
    "T22" = 0.1263194163093812
    "G23" = 67.91378400000001
    "H23" = 203.741352
    "I23" = 0.372191455479308
    "J23" = 0.3777393062881658
    "M23" = 5.85326
    "N23" = 17.55978
    "O23" = 0.1200026410479322
    "P23" = 0.1246091315254933
    "Q23" = 397.51703533584
    "R23" = 3577.65331802256
    "S23" = 0.04466395763299085
    "T23" = 0.04706976689961067
    "G24" = 67.91378400000001
    "H24" = 203.741352
    "I24" = 0.372191455479308
    "J24" = 0.3777393062881658
    "M24" = 10.959131
    "N24" = 32.877393
    "O24" = 0.2246824271585863
    "P24" = 0.2333072161810874
    "Q24" = 744.2760555617042
    "R24" = 6698.484500055337
    "S24" = 0.08362487958477785
    "T24" = 0.08812930599226708
    "G25" = 67.91378400000001
    "H25" = 203.741352
    "I25" = 0.372191455479308
    "J25" = 0.3777393062881658
    "M25" = 5.4093935
    "N25" = 10.818787
    "O25" = 0.1109025579706895
    "P25" = 0.07677315161290731
    "Q25" = 367.3723817300041
    "R25" = 2204.234290380024
    "S25" = 0.04127698446748925
    "T25" = 0.02900023703181578
    "G26" = 67.91378400000001
    "H26" = 203.741352
    "I26" = 0.372191455479308
    "J26" = 0.3777393062881658
    "M26" = 10.84612833333333
    "N26" = 32.538385
    "O26" = 0.2223656637745133
    "P26" = 0.230901520183746
    "Q26" = 736.60161686628
    "R26" = 6629.41455179652
    "S26" = 0.08276260004885853
    "T26" = 0.08722058005509113
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
